$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.719.60'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.827.64'
$ws.Range("E3").Value = '  +1.89%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.21'
$ws.Range("E5").Value = '  -0.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.68'
$ws.Range("E6").Value = '  +4.43%  '

$ws.Range("E7").Value = '  +1.68%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  +3.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.09'
$ws.Range("E10").Value = '  +1.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0847'
$ws.Range("E12").Value = '  +1.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.91'
$ws.Range("E13").Value = '  -0.72%  '

$ws.Range("E14").Value = '  +1.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.269.67'
$ws.Range("E15").Value = '  +1.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.982'
$ws.Range("E16").Value = '  +6.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.817.41'
$ws.Range("E17").Value = '  +1.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.738.19'
$ws.Range("E18").Value = '  +0.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.45'
$ws.Range("E19").Value = '  +11.58%  '

$ws.Range("E21").Value = '  +1.89%  '

$ws.Range("E22").Value = '  +1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.37'
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.86'
$ws.Range("E24").Value = '  +1.26%  '

$ws.Range("E25").Value = '  +1.29%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.23'
$ws.Range("E26").Value = '  +0.72%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.163'
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.55'
$ws.Range("E29").Value = '  +3.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.45'
$ws.Range("E30").Value = '  +6.42%  '

$ws.Range("E31").Value = '  +2.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.35'
$ws.Range("E32").Value = '  +1.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '52.80'
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0895'
$ws.Range("E34").Value = '  +8.33%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.62'
$ws.Range("E35").Value = '  +1.12%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0448'
$ws.Range("E36").Value = '  -1.13%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.93'
$ws.Range("E38").Value = '  +2.68%  '

$ws.Range("E39").Value = '  +1.93%  '

$ws.Range("E40").Value = '  +2.50%  '

$ws.Range("E41").Value = '  +1.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.52'
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.37'
$ws.Range("E43").Value = '  +2.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.21'
$ws.Range("E44").Value = '  +1.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.04'
$ws.Range("E45").Value = '  -0.35%  '

$ws.Range("E46").Value = '  +8.76%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.49'
$ws.Range("E47").Value = '  +7.03%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.168.33'
$ws.Range("E48").Value = '  +2.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.245'
$ws.Range("E49").Value = '  +22.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.953'
$ws.Range("E50").Value = '  +5.50%  '

$ws.Range("E51").Value = '  +1.81%  '
